# Commit: "fix for importer and more tests"
# - Teams sheet: selection collapses to the single cell A17 (was A2:I2 + A17)
# - Athletes sheet: add a new test value in B2 (new shared string "test"),
#   which extends the sheet's used range/dimension, and leaves the
#   selection on B2.

$wb = $excel.ActiveWorkbook

$teams = $wb.Worksheets.Item("Teams")
$teams.Range("A17").Select()

$athletes = $wb.Worksheets.Item("Athletes")
$athletes.Range("B2").Value = "test"
$athletes.Range("B2").Select()
